$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    # Force the cell to remain a text cell with the exact string value,
    # even when the string looks like a number (e.g. "240.42"), without
    # leaving a residual number-format style on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.319.86'
Set-TextValue $ws.Range("E2") '  -0.36%  '
Set-TextValue $ws.Range("D3") '1.842.22'
Set-TextValue $ws.Range("E3") '  -0.44%  '
Set-TextValue $ws.Range("D4") '0.9984'
Set-TextValue $ws.Range("E4") '  -0.23%  '
Set-TextValue $ws.Range("D5") '240.42'
Set-TextValue $ws.Range("D6") '0.6282'
Set-TextValue $ws.Range("E6") '  -0.15%  '
Set-TextValue $ws.Range("D7") '0.9990'
Set-TextValue $ws.Range("E7") '  -0.20%  '
Set-TextValue $ws.Range("D8") '0.07457'
Set-TextValue $ws.Range("E8") '  -2.50%  '
Set-TextValue $ws.Range("D9") '0.2894'
Set-TextValue $ws.Range("E9") '  -0.62%  '
Set-TextValue $ws.Range("D10") '24.33'
Set-TextValue $ws.Range("E10") '  -2.29%  '
Set-TextValue $ws.Range("D11") '0.07728'
Set-TextValue $ws.Range("E11") '  -0.13%  '
Set-TextValue $ws.Range("D12") '1.843.54'
Set-TextValue $ws.Range("E12") '  -2.37%  '
Set-TextValue $ws.Range("D13") '4.996'
Set-TextValue $ws.Range("E13") '  -0.75%  '
Set-TextValue $ws.Range("D14") '0.6774'
Set-TextValue $ws.Range("E14") '  -0.53%  '
Set-TextValue $ws.Range("D15") '0.00001011'
Set-TextValue $ws.Range("E15") '  -4.94%  '
Set-TextValue $ws.Range("D16") '81.97'
Set-TextValue $ws.Range("E16") '  -1.77%  '
Set-TextValue $ws.Range("D17") '6.121'
Set-TextValue $ws.Range("E17") '  -1.01%  '
Set-TextValue $ws.Range("D18") '29.364.58'
Set-TextValue $ws.Range("E18") '  -0.42%  '
Set-TextValue $ws.Range("D19") '228.32'
Set-TextValue $ws.Range("E19") '  -0.09%  '
Set-TextValue $ws.Range("D20") '12.29'
Set-TextValue $ws.Range("E20") '  -0.31%  '
Set-TextValue $ws.Range("E21") '  -0.24%  '
Set-TextValue $ws.Range("D22") '7.427'
Set-TextValue $ws.Range("E22") '  -0.51%  '
Set-TextValue $ws.Range("D23") '0.9990'
Set-TextValue $ws.Range("E23") '  -0.19%  '
Set-TextValue $ws.Range("D24") '158.90'
Set-TextValue $ws.Range("E24") '  +0.85%  '
Set-TextValue $ws.Range("D25") '0.1373'
Set-TextValue $ws.Range("E25") '  -0.90%  '
Set-TextValue $ws.Range("D26") '8.407'
Set-TextValue $ws.Range("E26") '  -0.32%  '
Set-TextValue $ws.Range("D27") '17.54'
Set-TextValue $ws.Range("E27") '  -0.84%  '
Set-TextValue $ws.Range("D28") '0.06502'
Set-TextValue $ws.Range("E28") '  +15.87%  '
Set-TextValue $ws.Range("D29") '1.391'
Set-TextValue $ws.Range("E29") '  +0.29%  '
Set-TextValue $ws.Range("D30") '1.472'
Set-TextValue $ws.Range("E30") '  +0.69%  '
Set-TextValue $ws.Range("E31") '  -1.21%  '
Set-TextValue $ws.Range("D32") '4.067'
Set-TextValue $ws.Range("E32") '  +0.37%  '
Set-TextValue $ws.Range("D33") '1.816'
Set-TextValue $ws.Range("E33") '  -1.59%  '
Set-TextValue $ws.Range("D34") '1.140'
Set-TextValue $ws.Range("E34") '  -2.09%  '
Set-TextValue $ws.Range("D35") '0.7032'
Set-TextValue $ws.Range("E35") '  +0.22%  '
Set-TextValue $ws.Range("D36") '2.580'
Set-TextValue $ws.Range("E36") '  -0.43%  '
Set-TextValue $ws.Range("D37") '1.258.56'
Set-TextValue $ws.Range("E37") '  +2.34%  '
Set-TextValue $ws.Range("D38") '2.830'
Set-TextValue $ws.Range("E38") '  +4.61%  '
Set-TextValue $ws.Range("E39") '  +0.34%  '
Set-TextValue $ws.Range("D40") '6.516'
Set-TextValue $ws.Range("E40") '  +1.15%  '
Set-TextValue $ws.Range("D41") '0.9108'
Set-TextValue $ws.Range("E41") '  +0.40%  '
Set-TextValue $ws.Range("D42") '0.9982'
Set-TextValue $ws.Range("E42") '  -0.29%  '
Set-TextValue $ws.Range("D43") '2.004.44'
Set-TextValue $ws.Range("E43") '  -10.72%  '
Set-TextValue $ws.Range("D44") '101.21'
Set-TextValue $ws.Range("E44") '  -1.10%  '
Set-TextValue $ws.Range("D45") '66.33'
Set-TextValue $ws.Range("E46") '  +1.09%  '
Set-TextValue $ws.Range("D47") '7.012'
Set-TextValue $ws.Range("E47") '  -2.54%  '
Set-TextValue $ws.Range("B48") 'EnergySwap'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D48") '9.032'
Set-TextValue $ws.Range("E48") '  +0.33%  '
Set-TextValue $ws.Range("B49") 'BabyDogeCoin'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D49") '0.00000000115'
Set-TextValue $ws.Range("E49") '  -4.30%  '
Set-TextValue $ws.Range("D50") '0.3943'
Set-TextValue $ws.Range("E50") '  -2.04%  '
Set-TextValue $ws.Range("E51") '  -0.63%  '
